$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.826.21"
$ws.Range("E2").Value = "  -0.26%  "

$ws.Range("D3").Value = "3.806.11"
$ws.Range("E3").Value = "  -1.40%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "705.60"
$ws.Range("E5").Value = "  +0.98%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.21"
$ws.Range("E6").Value = "  -1.97%  "

$ws.Range("D7").Value = "3.805.10"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.521"
$ws.Range("E9").Value = "  -0.63%  "

$ws.Range("E10").Value = "  -1.33%  "

$ws.Range("E11").Value = "  +5.45%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.460"
$ws.Range("E12").Value = "  -0.23%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000252"
$ws.Range("E13").Value = "  -1.85%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.98"
$ws.Range("E14").Value = "  -1.52%  "

$ws.Range("D15").Value = "4.448.34"
$ws.Range("E15").Value = "  -1.39%  "

$ws.Range("D16").Value = "3.835.75"
$ws.Range("E16").Value = "  -0.96%  "

$ws.Range("D17").Value = "70.805.52"
$ws.Range("E17").Value = "  -0.40%  "

$ws.Range("E18").Value = "  +0.03%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.13"
$ws.Range("E19").Value = "  -1.56%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.28"
$ws.Range("E20").Value = "  -2.55%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "496.19"
$ws.Range("E21").Value = "  +0.02%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.66"
$ws.Range("E22").Value = "  -4.62%  "

$ws.Range("E23").Value = "  +0.18%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.94"
$ws.Range("E24").Value = "  -1.40%  "

$ws.Range("E25").Value = "  -0.73%  "

$ws.Range("B26").Value = "WrappedeETH"
$ws.Range("C26").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D26").Value = "3.957.92"
$ws.Range("E26").Value = "  -1.17%  "

$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.09"
$ws.Range("E27").Value = "  -1.92%  "

$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.38"
$ws.Range("E28").Value = "  -3.23%  "

$ws.Range("E29").Value = "  +0.10%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.04"
$ws.Range("E30").Value = "  -5.34%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.08"
$ws.Range("E31").Value = "  -3.36%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.27"
$ws.Range("E32").Value = "  -0.46%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.32"
$ws.Range("E33").Value = "  -4.38%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.03"
$ws.Range("E34").Value = "  -2.14%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.175"
$ws.Range("E35").Value = "  -4.03%  "

$ws.Range("D36").Value = "3.775.58"
$ws.Range("E36").Value = "  -0.96%  "

$ws.Range("B37").Value = "Binance-PegBSC-USD"
$ws.Range("C37").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  -0.03%  "

$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "9.09"
$ws.Range("E38").Value = "  -2.11%  "

$ws.Range("E39").Value = "  -2.91%  "

$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.33"
$ws.Range("E40").Value = "  -2.68%  "

$ws.Range("B41").Value = "Mantle"
$ws.Range("C41").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.03"
$ws.Range("E41").Value = "  -0.35%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.95"
$ws.Range("E42").Value = "  -1.69%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.24"
$ws.Range("E43").Value = "  -4.28%  "

$ws.Range("E44").Value = "  -0.01%  "

$ws.Range("B45").Value = "FLOKI"
$ws.Range("C45").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.000324"
$ws.Range("E45").Value = "  +4.42%  "

$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").Value = "  +0.13%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "165.30"
$ws.Range("E47").Value = "  +1.44%  "

$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "426.46"
$ws.Range("E48").Value = "  +2.10%  "

$ws.Range("B49").Value = "OKB"
$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "48.82"
$ws.Range("E49").Value = "  +0.29%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.61"
$ws.Range("E50").Value = "  -1.18%  "

$ws.Range("B51").Value = "ONDO"
$ws.Range("C51").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.36"
$ws.Range("E51").Value = "  -1.76%  "

